# Update country case-count table and the "last updated" timestamp
# to match the refreshed data snapshot (24 Mar 2020, 07:46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last-updated timestamp (row 1)
$ws.Range('A1').Value = 'Datos actualizados a 24 de Marzo de 2020 a las 07:46'

# Row 18
$ws.Range('A18').Value = 'Australia'
$ws.Range('B18').Value = 2136
$ws.Range('C18').Value = 249
$ws.Range('D18').Value = 118
$ws.Range('E18').Value = 2010
$ws.Range('F18').Value = 11
$ws.Range('G18').Value = 1
$ws.Range('H18').Value = 8

# Row 19
$ws.Range('A19').Value = 'Canada'
$ws.Range('B19').Value = 2091
$ws.Range('D19').Value = 320
$ws.Range('E19').Value = 1747
$ws.Range('F19').Value = 1
$ws.Range('H19').Value = 24

# Row 20
$ws.Range('A20').Value = 'Portugal'
$ws.Range('B20').Value = 2060
$ws.Range('D20').Value = 14
$ws.Range('E20').Value = 2023
$ws.Range('F20').Value = 47
$ws.Range('H20').Value = 23

# Row 21
$ws.Range('A21').Value = 'Suecia'
$ws.Range('B21').Value = 2046
$ws.Range('C21').Value = 0
$ws.Range('D21').Value = 16
$ws.Range('E21').Value = 2003
$ws.Range('F21').Value = 104
$ws.Range('G21').Value = 0
$ws.Range('H21').Value = 27

# Row 31
$ws.Range('A31').Value = 'Pakistan'
$ws.Range('B31').Value = 892
$ws.Range('C31').Value = 17
$ws.Range('D31').Value = 13
$ws.Range('E31').Value = 873
$ws.Range('F31').Value = 0
$ws.Range('H31').Value = 6

# Row 32
$ws.Range('A32').Value = 'Luxemburgo'
$ws.Range('D32').Value = 6
$ws.Range('E32').Value = 861
$ws.Range('F32').Value = 3
$ws.Range('H32').Value = 8

# Row 66
$ws.Range('A66').Value = 'Taiwan'
$ws.Range('B66').Value = 215
$ws.Range('C66').Value = 20
$ws.Range('D66').Value = 29
$ws.Range('E66').Value = 184
$ws.Range('F66').Value = 0
$ws.Range('H66').Value = 2

# Row 67
$ws.Range('A67').Value = 'Bulgaria'
$ws.Range('B67').Value = 201
$ws.Range('D67').Value = 3
$ws.Range('E67').Value = 195
$ws.Range('F67').Value = 8
$ws.Range('H67').Value = 3

# Row 68
$ws.Range('A68').Value = 'Emiratos Arabes Unidos'
$ws.Range('B68').Value = 198
$ws.Range('D68').Value = 41
$ws.Range('E68').Value = 155
$ws.Range('F68').Value = 2

# Row 71
$ws.Range('A71').Value = 'Hungria'
$ws.Range('B71').Value = 187
$ws.Range('C71').Value = 20
$ws.Range('D71').Value = 21
$ws.Range('E71').Value = 158
$ws.Range('F71').Value = 6
$ws.Range('H71').Value = 8

# Row 72
$ws.Range('A72').Value = 'Eslovaquia'
$ws.Range('B72').Value = 186
$ws.Range('D72').Value = 7
$ws.Range('F72').Value = 2

# Row 73
$ws.Range('A73').Value = 'Letonia'
$ws.Range('B73').Value = 180
$ws.Range('E73').Value = 179
$ws.Range('F73').Value = 0
$ws.Range('H73').Value = 0

# Row 74
$ws.Range('A74').Value = 'Lituania'
$ws.Range('B74').Value = 179
$ws.Range('D74').Value = 1
$ws.Range('E74').Value = 177
$ws.Range('F74').Value = 1
$ws.Range('H74').Value = 1

# Row 95
$ws.Range('A95').Value = 'Oman'
$ws.Range('B95').Value = 84
$ws.Range('C95').Value = 18
$ws.Range('D95').Value = 17
$ws.Range('E95').Value = 67

# Row 96
$ws.Range('A96').Value = 'Bielorrusia'
$ws.Range('B96').Value = 81
$ws.Range('D96').Value = 22
$ws.Range('E96').Value = 59

# Row 97
$ws.Range('A97').Value = 'Senegal'
$ws.Range('B97').Value = 79
$ws.Range('D97').Value = 8
$ws.Range('E97').Value = 71
$ws.Range('H97').Value = 0

# Row 98
$ws.Range('A98').Value = 'Ucrania'
$ws.Range('B98').Value = 73
$ws.Range('D98').Value = 1
$ws.Range('E98').Value = 69
$ws.Range('F98').Value = 0
$ws.Range('H98').Value = 3

# Row 99
$ws.Range('A99').Value = 'Azerbaiyan'
$ws.Range('B99').Value = 72
$ws.Range('D99').Value = 10
$ws.Range('E99').Value = 61
$ws.Range('F99').Value = 3
$ws.Range('H99').Value = 1

# Row 100
$ws.Range('A100').Value = 'Reunion'
$ws.Range('B100').Value = 71
$ws.Range('C100').Value = 0
$ws.Range('D100').Value = 1
$ws.Range('E100').Value = 70

# Row 101
$ws.Range('A101').Value = 'Kazajistan'
$ws.Range('B101').Value = 68
$ws.Range('C101').Value = 6
$ws.Range('D101').Value = 0
$ws.Range('E101').Value = 68

# Row 102
$ws.Range('A102').Value = 'Georgia'
$ws.Range('B102').Value = 66
$ws.Range('C102').Value = 5
$ws.Range('D102').Value = 9
$ws.Range('E102').Value = 57
$ws.Range('F102').Value = 1
$ws.Range('H102').Value = 0

# Row 103
$ws.Range('A103').Value = 'Guadalupe'
$ws.Range('B103').Value = 62
$ws.Range('D103').Value = 0
$ws.Range('E103').Value = 61
$ws.Range('F103').Value = 4
$ws.Range('H103').Value = 1
